# "excel read/write enum implementation"
#
# Inserts a new "EnumTaskTest" worksheet (a copy of the "Tasks" sheet's
# name/last name/age/Amount/description/date columns, plus a new "Status"
# enum column) right before the existing "DateTask" sheet, so the tab
# order becomes: Tasks, EnumTaskTest, DateTask, People, TaskMiddle.

$wb = $excel.ActiveWorkbook

$tasks = $wb.Worksheets.Item("Tasks")

# Worksheets.Add() with no arguments inserts the new sheet immediately
# before the currently active sheet. The workbook's active sheet is
# "DateTask" (activeTab=1), so this lands the new tab exactly between
# "Tasks" and "DateTask".
$enumTaskTest = $wb.Worksheets.Add()
$enumTaskTest.Name = "EnumTaskTest"

# Re-resolve "DateTask" by name now that the sheet collection has shifted.
$dateTask = $wb.Worksheets.Item("DateTask")

# Clone the shared name/last name/age/Amount/description/date block
# (values + number formats/styles) straight from "Tasks".
$tasks.Range("A1:F11").Copy($enumTaskTest.Range("A1:F11"))

# New "Status" enum column.
$enumTaskTest.Range("G1").Value = "Status"
$enumTaskTest.Range("G2").Value = "ACTIVE"
$enumTaskTest.Range("G3").Value = "FINISHED"
$enumTaskTest.Range("G4").Value = "FINISHED"
$enumTaskTest.Range("G5").Value = "FINISHED"
$enumTaskTest.Range("G6").Value = "FINISHED"
$enumTaskTest.Range("G7").Value = "FINISHED"
$enumTaskTest.Range("G8").Value = "FINISHED"
$enumTaskTest.Range("G9").Value = "FINISHED"
$enumTaskTest.Range("G10").Value = "ACTIVE"
$enumTaskTest.Range("G11").Value = "ACTIVE"

$enumTaskTest.Columns.Item(6).ColumnWidth = 14.998697916666666

# Restore "DateTask"'s remembered selection, then make "EnumTaskTest" the
# active tab with its own remembered selection (selecting it last is what
# leaves it as the active/visible sheet).
$dateTask.Range("G1").Select()
$enumTaskTest.Range("H8").Select()
